# Update of league bases, swapping the two mis-ordered match rows in each
# of the three pairs below (columns B..AD swap while column A, the row's
# running index, stays put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 30   # column AD

function Swap-Rows($rowA, $rowB) {
    $rangeA = $ws.Range($ws.Cells.Item($rowA, 2), $ws.Cells.Item($rowA, $lastCol))
    $rangeB = $ws.Range($ws.Cells.Item($rowB, 2), $ws.Cells.Item($rowB, $lastCol))

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}

Swap-Rows 120 121
Swap-Rows 158 159
Swap-Rows 171 172
